$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Industries" column (H) policy flag is turned off (1 -> 0)
# for rows 31 through 176.
$ws.Range("H31:H176").Value = 0
